# Update "想去人数" (F column) counts across the sheets to the values
# recorded in the new data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 998
$ws.Range("F4").Value = 239
$ws.Range("F6").Value = 1163
$ws.Range("F7").Value = 944
$ws.Range("F9").Value = 65
$ws.Range("F11").Value = 904
$ws.Range("F12").Value = 326
$ws.Range("F14").Value = 530
$ws.Range("F15").Value = 1379
$ws.Range("F17").Value = 1283
$ws.Range("F18").Value = 2943
$ws.Range("F19").Value = 313
$ws.Range("F20").Value = 1567
$ws.Range("F21").Value = 1318
$ws.Range("F22").Value = 758
$ws.Range("F23").Value = 218
$ws.Range("F26").Value = 1078
$ws.Range("F28").Value = 3335
$ws.Range("F29").Value = 650
$ws.Range("F31").Value = 1475

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 12
$ws.Range("F9").Value = 39
$ws.Range("F14").Value = 2

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 778

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 778
$ws.Range("F6").Value = 998
$ws.Range("F7").Value = 239
$ws.Range("F10").Value = 1163
$ws.Range("F11").Value = 944
$ws.Range("F14").Value = 65
$ws.Range("F18").Value = 12
$ws.Range("F19").Value = 39
$ws.Range("F23").Value = 904
$ws.Range("F24").Value = 326
$ws.Range("F26").Value = 530
$ws.Range("F27").Value = 1379
$ws.Range("F29").Value = 1283
$ws.Range("F30").Value = 2943
$ws.Range("F31").Value = 313
$ws.Range("F32").Value = 1567
$ws.Range("F33").Value = 1318
$ws.Range("F34").Value = 758
$ws.Range("F35").Value = 218
$ws.Range("F40").Value = 1078
$ws.Range("F42").Value = 3335
$ws.Range("F43").Value = 650
$ws.Range("F45").Value = 1475
$ws.Range("F46").Value = 2
